# Seed database update: populate "Machine Group Name" (column C) for every
# machine row, and add two new lookup strings used by that column.
#
# Rows 2-16  -> "SMS_Group_1, TCM_Group_1"
# Rows 17-66 -> "SMS_Group_1"
#
# Write the longer/second-appearing string group first so the shared-string
# table receives the two new entries in the same order Excel produced them
# (index 105 = "SMS_Group_1", index 106 = "SMS_Group_1, TCM_Group_1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17:C66").Value = "SMS_Group_1"
$ws.Range("C2:C16").Value = "SMS_Group_1, TCM_Group_1"

# Column C now holds longer text than before, so widen it to fit (matches
# the author's manual resize of column C from ~20.8 to ~26.2 characters).
$ws.Columns.Item(3).ColumnWidth = 25.333333333333332

# Reflect the author's final scroll position / cell selection.
$ws.Range("B20").Select() | Out-Null
